$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Weekly")

$ws.Range("D2").Value = "Entity Code"
$ws.Range("D3").Value = "Entity Name"

$ws.Activate() | Out-Null
$ws.Range("D4").Select() | Out-Null
